$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.131.33"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.772.08"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "628.16"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").Value = "165.91"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("D7").Value = "3.767.95"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").Value = "6.77"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "34.95"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "4.408.38"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "3.774.04"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "69.144.01"
$ws.Range("E18").Value = "  -3.11%  "
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "7.00"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").Value = "463.37"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "9.53"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").Value = "0.707"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "82.94"
$ws.Range("D25").Value = "0.0000144"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").Value = "11.97"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "2.14"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "3.922.87"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").Value = "2.27"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").Value = "7.10"
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").Value = "28.55"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").Value = "0.173"
$ws.Range("E35").Value = "  +15.80%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "3.727.49"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "8.99"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "3.30"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "0.958"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "157.56"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").Value = "1.42"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("E47").Value = "  +3.80%  "
$ws.Range("D48").Value = "43.09"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "46.66"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "8.36"
$ws.Range("E51").Value = "  -0.50%  "
